# PLUS DE BUG !!!!!!!!!!!!!!!! CORRECTION TOUT ÇA
#
# This script reproduces two kinds of changes that were made to the
# "S3BFI.xlsx" planning workbook:
#
#   1. A handful of worksheet tab names that exceeded Excel's 31
#      character limit for sheet names got truncated to 31 characters.
#
#   2. The "Organisation détaillée" weekly-session tables (columns
#      A..D, rows 35 and below) on every "R3.xx" resource sheet had
#      their session date / Amphi / TD / TP cells cleared out.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Truncate worksheet names to Excel's 31-character tab-name limit
# ---------------------------------------------------------------
$renames = @{
    "R3.07 SQL dans un langage de programmation"   = "R3.07 SQL dans un langage de pr"
    "R3.10 Management des systèmes d'information " = "R3.10 Management des systèmes d"
    "R3.11 Droit des contrats et du numérique"      = "R3.11 Droit des contrats et du "
    "R3.13 Communication professionnelle"           = "R3.13 Communication professionn"
    "R3.B.L1 Compléments cryptographie"             = "R3.B.L1 Compléments cryptograph"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# ---------------------------------------------------------------
# 2) Clear out the planning cells (A:D, rows 35-end) on every
#    resource sheet
# ---------------------------------------------------------------
$clearRanges = @{
    "R3.01 Développement Web"                       = "A35:D39"
    "R3.02 Développement efficace"                  = "A35:D39"
    "R3.03 Analyse"                                 = "A35:D38"
    "R3.04 Qualité de développement"                = "A35:D41"
    "R3.05 Programmation système"                   = "A35:D41"
    "R3.06 Architecture des réseaux"                = "A35:D37"
    "R3.07 SQL dans un langage de pr"               = "A35:D38"
    "R3.08 Probabilités"                            = "A35:D43"
    "R3.09 Cryptographie et sécurité"                = "A35:D43"
    "R3.10 Management des systèmes d"               = "A35:D42"
    "R3.11 Droit des contrats et du "                = "A35:D41"
    "R3.12 Anglais"                                 = "A35:D44"
    "R3.13 Communication professionn"               = "A35:D40"
    "R3.14 PPP"                                     = "A35:D37"
    "R3.B.L1 Compléments cryptograph"               = "A35:D37"
}

foreach ($sheetName in $clearRanges.Keys) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range($clearRanges[$sheetName]).ClearContents()
}
